$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value would otherwise be auto-detected as a number
# by Excel need to be explicitly formatted as Text first, so the value is
# stored as a string (matching the source data which is always textual).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "56.363.07"
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").Value = "2.479.38"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").Value = "488.17"
$ws.Range("E5").Value = "  +4.10%  "
$ws.Range("D6").Value = "146.12"
$ws.Range("E6").Value = "  +8.63%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "0.513"
$ws.Range("E8").Value = "  +4.22%  "
$ws.Range("D9").Value = "2.497.28"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "5.76"
$ws.Range("E10").Value = "  +7.66%  "
$ws.Range("D11").Value = "0.0972"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("E12").Value = "  +3.58%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "2.931.31"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "56.396.60"
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").Value = "21.23"
$ws.Range("E16").Value = "  +6.07%  "
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("D18").Value = "2.496.96"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("E19").Value = "  +7.03%  "
$ws.Range("D20").Value = "10.16"
$ws.Range("E20").Value = "  +6.95%  "
$ws.Range("D21").Value = "319.74"
$ws.Range("E21").Value = "  +2.65%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +7.85%  "
$ws.Range("D24").Value = "58.63"
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("E25").Value = "  +5.94%  "
$ws.Range("E26").Value = "  +6.23%  "
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").Value = "2.603.68"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").Value = "7.74"
$ws.Range("E29").Value = "  +6.16%  "
$ws.Range("D30").Value = "0.0₃0791"
$ws.Range("E30").Value = "  +7.92%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "149.13"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("D36").Value = "1.15"
$ws.Range("E36").Value = "  +6.62%  "
$ws.Range("D37").Value = "3.73"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("E38").Value = "  +6.60%  "
$ws.Range("D39").Value = "34.12"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("E40").Value = "  +6.23%  "
$ws.Range("D41").Value = "0.615"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("E42").Value = "  +4.63%  "
$ws.Range("D43").Value = "0.994"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("E44").Value = "  +6.11%  "
$ws.Range("E45").Value = "  +11.21%  "
$ws.Range("D46").Value = "259.14"
$ws.Range("E46").Value = "  +15.56%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("D49").Value = "0.0913"
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("D50").Value = "1.900.83"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").Value = "17.62"
$ws.Range("E51").Value = "  +4.70%  "
